# Scheduled market-data refresh: update currentAveragePrice / LevePrice / LeveProfit
# columns (H-N) across the Leve profit sheets with freshly pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 3724.8572
$ws.Range("I11").Value = 3724.8572
$ws.Range("K11").Value = 3724.8572
$ws.Range("M11").Value = -3584.8572

# Row 123
$ws.Range("H123").Value = 39273.332
$ws.Range("J123").Value = 39273.332
$ws.Range("L123").Value = 39273.332
$ws.Range("N123").Value = -49073.332

# Row 129
$ws.Range("H129").Value = 880.4167
$ws.Range("J129").Value = 1031.3334
$ws.Range("L129").Value = 3094.0002
$ws.Range("N129").Value = -13094.0002

# Row 132
$ws.Range("H132").Value = 235013.67
$ws.Range("I132").Value = 235013.67
$ws.Range("K132").Value = 705041.01
$ws.Range("M132").Value = -702511.01

# Row 138
$ws.Range("H138").Value = 2203.879
$ws.Range("I138").Value = 1688.738
$ws.Range("J138").Value = 3105.375
$ws.Range("K138").Value = 5066.214
$ws.Range("L138").Value = 9316.125
$ws.Range("M138").Value = 73.78600000000006
$ws.Range("N138").Value = -19596.125

# Row 141
$ws.Range("H141").Value = 1249.2354
$ws.Range("I141").Value = 682.775
$ws.Range("J141").Value = 3309.0908
$ws.Range("K141").Value = 2048.325
$ws.Range("L141").Value = 9927.2724
$ws.Range("M141").Value = 3131.675
$ws.Range("N141").Value = -20287.2724

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 795868.7
$ws.Range("I2").Value = 685.12
$ws.Range("J2").Value = 2452501.2
$ws.Range("K2").Value = 685.12
$ws.Range("L2").Value = 2452501.2
$ws.Range("M2").Value = -572.12
$ws.Range("N2").Value = -2452727.2

# Row 32
$ws.Range("H32").Value = 2865.44
$ws.Range("I32").Value = 2639.694
$ws.Range("K32").Value = 2639.694
$ws.Range("M32").Value = -2352.694

# Row 116
$ws.Range("H116").Value = 795868.7
$ws.Range("I116").Value = 685.12
$ws.Range("J116").Value = 2452501.2
$ws.Range("K116").Value = 685.12
$ws.Range("L116").Value = 2452501.2
$ws.Range("M116").Value = 1608.88
$ws.Range("N116").Value = -2457089.2

# Row 132
$ws.Range("H132").Value = 772884.25
$ws.Range("I132").Value = 822627.0600000001
$ws.Range("J132").Value = 250585
$ws.Range("K132").Value = 2467881.18
$ws.Range("L132").Value = 751755
$ws.Range("M132").Value = -2465351.18
$ws.Range("N132").Value = -756815

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 795868.7
$ws.Range("I3").Value = 685.12
$ws.Range("J3").Value = 2452501.2
$ws.Range("K3").Value = 685.12
$ws.Range("L3").Value = 2452501.2
$ws.Range("M3").Value = -571.12
$ws.Range("N3").Value = -2452729.2

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 382452.5
$ws.Range("I31").Value = 1186.4445
$ws.Range("J31").Value = 897161.7
$ws.Range("K31").Value = 1186.4445
$ws.Range("L31").Value = 897161.7
$ws.Range("M31").Value = -891.4445000000001
$ws.Range("N31").Value = -897751.7

# Row 34
$ws.Range("H34").Value = 382452.5
$ws.Range("I34").Value = 1186.4445
$ws.Range("J34").Value = 897161.7
$ws.Range("K34").Value = 1186.4445
$ws.Range("L34").Value = 897161.7
$ws.Range("M34").Value = -984.4445000000001
$ws.Range("N34").Value = -897565.7

# Row 58
$ws.Range("H58").Value = 2445.06
$ws.Range("I58").Value = 1167.5416
$ws.Range("J58").Value = 3624.3076
$ws.Range("K58").Value = 1167.5416
$ws.Range("L58").Value = 3624.3076
$ws.Range("M58").Value = -964.5416
$ws.Range("N58").Value = -4030.3076

# Row 132
$ws.Range("H132").Value = 2784
$ws.Range("I132").Value = 2741.524
$ws.Range("J132").Value = 3007
$ws.Range("K132").Value = 8224.572
$ws.Range("L132").Value = 9021
$ws.Range("M132").Value = -5694.572
$ws.Range("N132").Value = -14081

# Row 134
$ws.Range("H134").Value = 2199.524
$ws.Range("I134").Value = 2402.889
$ws.Range("J134").Value = 979.3333
$ws.Range("K134").Value = 7208.667
$ws.Range("L134").Value = 2937.9999
$ws.Range("M134").Value = -4673.667
$ws.Range("N134").Value = -8007.9999

# Row 136
$ws.Range("H136").Value = 2445.06
$ws.Range("I136").Value = 1167.5416
$ws.Range("J136").Value = 3624.3076
$ws.Range("K136").Value = 3502.6248
$ws.Range("L136").Value = 10872.9228
$ws.Range("M136").Value = -952.6248000000001
$ws.Range("N136").Value = -15972.9228

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 8929322
$ws.Range("I5").Value = 364.82352
$ws.Range("J5").Value = 12821431
$ws.Range("K5").Value = 1094.47056
$ws.Range("L5").Value = 38464293
$ws.Range("M5").Value = -982.47056
$ws.Range("N5").Value = -38464517

# Row 131
$ws.Range("H131").Value = 3275.318
$ws.Range("J131").Value = 2495.2903
$ws.Range("L131").Value = 7485.8709
$ws.Range("N131").Value = -17565.8709

# Row 135
$ws.Range("H135").Value = 8929322
$ws.Range("I135").Value = 364.82352
$ws.Range("J135").Value = 12821431
$ws.Range("K135").Value = 3283.41168
$ws.Range("L135").Value = 115392879
$ws.Range("M135").Value = -748.4116799999997
$ws.Range("N135").Value = -115397949

$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 2305.5
$ws.Range("I126").Value = 2305.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6916.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -4446.5
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 1777.5714
$ws.Range("I40").Value = 1777.5714
$ws.Range("K40").Value = 1777.5714
$ws.Range("M40").Value = -1641.5714

# Row 132
$ws.Range("H132").Value = 3851.225
$ws.Range("I132").Value = 4208
$ws.Range("K132").Value = 12624
$ws.Range("M132").Value = -10094

# Row 136
$ws.Range("H136").Value = 1731.3784
$ws.Range("I136").Value = 1118.5555
$ws.Range("J136").Value = 3386
$ws.Range("K136").Value = 3355.6665
$ws.Range("L136").Value = 10158
$ws.Range("M136").Value = -805.6664999999998
$ws.Range("N136").Value = -15258

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 14471111
$ws.Range("I2").Value = 2891428.5
$ws.Range("J2").Value = 55000000
$ws.Range("K2").Value = 2891428.5
$ws.Range("L2").Value = 55000000
$ws.Range("M2").Value = -2891316.5
$ws.Range("N2").Value = -55000224

# Row 124
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

# Row 132
$ws.Range("H132").Value = 3512.6155
$ws.Range("I132").Value = 4103.0234
$ws.Range("J132").Value = 691.7778
$ws.Range("K132").Value = 12309.0702
$ws.Range("L132").Value = 2075.3334
$ws.Range("M132").Value = -9779.0702
$ws.Range("N132").Value = -7135.3334

# Row 136
$ws.Range("H136").Value = 4427.82
$ws.Range("I136").Value = 4782.6523
$ws.Range("J136").Value = 347.25
$ws.Range("K136").Value = 14347.9569
$ws.Range("L136").Value = 1041.75
$ws.Range("M136").Value = -11797.9569
$ws.Range("N136").Value = -6141.75
